$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KnownIssues")

# Exercise number changed from 2 to 3
$ws.Range("D4").Value = 3

# Row 8: Type / Sub Type / Description / Severity
$ws.Range("C8").Value = "Code"
$ws.Range("D8").Value = "Function is too long"
$ws.Range("E8").Value = "function calculate - switch cases for 3 operators - can't reduse number of cases,because we want a variaty of questions with 3 operators"
$ws.Range("F8").Value = "Low"

# The old text used to live in the "Reason notes" column (K8); clear it
$ws.Range("K8").ClearContents()

# Update current selection to E8
$ws.Range("E8").Select()
